$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "fio" value in C2: replace the literal "\r\n" marker with a plain space.
$ws.Range("C2").Value = "Юрий Теуш"

# Move the active selection to C2 (was F4).
$ws.Range("C2").Select()
